$d = $word.ActiveDocument

# --- Comment 1 -------------------------------------------------------
# Anchors the whole tracked-change sentence describing how the line
# features were overlaid with the bathymetry grid to assign depths to
# the start/end points of each string.
$find1 = $d.Content.Find
$found1 = $find1.Execute(
    "We then overlaid these line features with a Using NGDC composite bathymetry grid developed by Feist et al. (2021) to provide assign a depth for to the start and end points of each stringeach point",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $rng1 = $find1.Parent.Duplicate
    $c1 = $d.Comments.Add($rng1, "The traps are simulated to the stringline first, and then bathymetry is extract to all points (not just the end coordinates)")
    $c1.Author = "Leena Riekkola"
    $c1.Initial = "LR"
}

# --- Comment 2 -------------------------------------------------------
# Anchors the single word "fully" inside the (re-worded) sentence about
# strings occurring fully in deep water.
$find2 = $d.Content.Find
$found2 = $find2.Execute("fully", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $rng2 = $find2.Parent.Duplicate
    $c2 = $d.Comments.Add($rng2, "Completely or partially")
    $c2.Author = "Leena Riekkola"
    $c2.Initial = "LR"
}
